# A new September transaction ("axis" @ 2024-09-11 06:57:42) was recorded
# at the top of the "2024" sheet's September log. This pushes every
# existing row at/after row 36 down by one row (the sheet's per-month
# "Details"/"Date" column pairs are append-at-top logs, so inserting a
# physical worksheet row reproduces the same downward shift seen for all
# of September/August/etc. below it, including the trailing "Broadband"
# label row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at 36, shifting rows 36:131 down to 37:132.
$ws.Rows.Item(36).Insert()

# Populate the newly inserted row with the latest September entry.
$ws.Range("R36").Value = "axis"
$ws.Range("S36").Value = "2024-09-11 06:57:42"
